# FIN13Final.xlsx — "Periodo and User Controller finally api.php"
# On sheet "Hoja1. Actividades":
#   - J10 ("se realizó") is renamed to "Activo"
#   - Row 11 is filled in with a new sample activity row.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1. Actividades")

$ws.Range("J10").Value = "Activo"

$ws.Range("A11").Value = "maicol es el mejor tipo que conozco"
$ws.Range("B11").Value = "maicol"
$ws.Range("C11").Value = "$$$"
$ws.Range("G11").Value = "X"
$ws.Range("H11").Value = "X"
$ws.Range("I11").Value = "billetera"
$ws.Range("J11").Value = "Activo"
